# Daily attendance processing - 2026-01-01 21:31:00
# Swap the "Recorded By" text for rows that were recorded by both
# dnasr281@gmail.com and the System, except for the sessions dated
# 23/12/2025 and 31/12/2025 (those stay untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"
$skipDates = @("23/12/2025", "31/12/2025")

for ($r = 2; $r -le $lastRow; $r++) {
    $recordedBy = $ws.Cells.Item($r, 7).Text
    if ($recordedBy -eq $oldValue) {
        $sessionDate = $ws.Cells.Item($r, 5).Text
        if (-not ($skipDates -contains $sessionDate)) {
            $ws.Cells.Item($r, 7).Value = $newValue
        }
    }
}
